$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 8.856799999999998
$ws.Range("D4").Value = -8.211000000000006
$ws.Range("D7").Value = -7.519400000000003
$ws.Range("A9").Value = -22.1751
$ws.Range("B9").Value = 6.569700000000007
$ws.Range("C9").Value = -12.05619999999999
$ws.Range("D11").Value = -6.91229999999999
$ws.Range("D15").Value = -8.130599999999992
$ws.Range("A18").Value = -22.14120000000001
$ws.Range("A20").Value = -19.92159999999998
$ws.Range("B23").Value = 8.907799999999998
$ws.Range("B24").Value = 4.959899999999999
$ws.Range("B26").Value = 5.625700000000003
$ws.Range("A27").Value = -21.95189999999999
$ws.Range("D30").Value = -7.248999999999995
$ws.Range("C32").Value = -12.48710000000001
$ws.Range("B34").Value = 9.508000000000008
$ws.Range("B35").Value = 8.70070000000001
$ws.Range("C38").Value = -11.9412
$ws.Range("D39").Value = -8.193799999999994
$ws.Range("D43").Value = -7.332700000000001
$ws.Range("C45").Value = -14.06169999999999
$ws.Range("D47").Value = -7.397799999999997
$ws.Range("B48").Value = 5.501300000000006
$ws.Range("C51").Value = -11.39549999999999
$ws.Range("B52").Value = 5.654499999999997
$ws.Range("C57").Value = -13.99649999999999
$ws.Range("C64").Value = -10.22339999999999
$ws.Range("B66").Value = 5.899699999999997
$ws.Range("B67").Value = 5.345800000000002
$ws.Range("A69").Value = -21.63429999999999
$ws.Range("D75").Value = -7.629499999999998
$ws.Range("A76").Value = -19.58219999999999
$ws.Range("B80").Value = 9.684299999999997
$ws.Range("A82").Value = -21.67580000000001
$ws.Range("D91").Value = -7.707199999999999
$ws.Range("D92").Value = -6.377600000000002
$ws.Range("C93").Value = -10.42139999999999
$ws.Range("B99").Value = 5.7197
